# Commit: "removed MORE hard coded strings......"
# Adds 16 new translation-key rows (190-205) to Sheet1, each holding an
# Android string-resource key in column A and its English value in column B
# (columns C-F, i.e. es/fr/pt/ru, are left blank for later translation -
# matching the source diff, which only populates A/B for these rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlVAlignCenter = -4108

# row, key, valignKey(bool), englishText, style ("plain" | "fillcenter")
$rows = @(
    @(190, "enter_enum_area_name",      $true,  "Enter the Enumeration Area name",            "plain"),
    @(191, "select_enumeration",        $true,  "Select an Enumeration",                      "fillcenter"),
    @(192, "select_time",               $true,  "Select Time",                                "plain"),
    @(193, "select_date",               $true,  "Select Date",                                "plain"),
    @(194, "your_pin_is",               $true,  "Your PIN is:",                               "plain"),
    @(195, "incorrect_answer_message",  $true,  "Oops! Incorrect answer.  Please try again.",  "plain"),
    @(196, "enter_other_question",      $true,  "Enter Other Question",                       "plain"),
    @(197, "reenter_pin",               $false, "Please re-enter the PIN",                    "plain"),
    @(198, "enter_pin",                 $true,  "Please enter a PIN",                         "plain"),
    @(199, "other_question",            $false, "Other Question",                             "plain"),
    @(200, "enter_answer",              $false, "Please enter an answer.",                    "plain"),
    @(201, "pin_not_match",             $false, "The PIN's do not match",                     "plain"),
    @(202, "enter_file_name",           $false, "You must enter a file name",                 "plain"),
    @(203, "pin_incorrect",             $false, "The current PIN is incorrect",               "plain"),
    @(204, "config_not_found",          $false, "Fatal! Config not found.",                   "plain"),
    @(205, "missing_parameter_rule",    $false, "Fatal! Missing required parameter: role",    "plain")
)

foreach ($row in $rows) {
    $r = $row[0]
    $key = $row[1]
    $keyCentered = $row[2]
    $text = $row[3]
    $style = $row[4]

    $keyCell = $ws.Cells.Item($r, 1)
    $keyCell.Value = $key
    if ($keyCentered) {
        $keyCell.VerticalAlignment = $xlVAlignCenter
    }

    $textCell = $ws.Cells.Item($r, 2)
    $textCell.Value = $text
    if ($style -eq "fillcenter") {
        $textCell.VerticalAlignment = $xlVAlignCenter
        $textCell.Interior.Color = 13431551
    }
}

# Match the author's final viewport/selection (best-effort - the cell that
# was active when the workbook was saved).
$ws.Range("B208").Select()
